$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine rows 2-4 into a single inline string in A2, then remove rows 3-4.
$ws.Range("A2").Value = "('Elemental Shaman', ['Token Creature — Elemental Shaman', '3/1'])"

$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()
